$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.849566698074341
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.746504545211792
$ws.Range("D1").Value = 1.276583552360535
$ws.Range("E1").Value = 0.9300251007080078
